$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ D = "27.787.59"; E = "  -2.07%  " }
    3 = @{ D = "1.753.55"; E = "  -4.08%  " }
    4 = @{ D = "1.008"; E = "  +0.45%  " }
    5 = @{ D = "321.05"; E = "  -2.72%  " }
    6 = @{ E = "  +0.27%  " }
    7 = @{ D = "0.4256"; E = "  -4.87%  " }
    8 = @{ D = "0.3639"; E = "  -3.88%  " }
    9 = @{ D = "43.24"; E = "  -3.35%  " }
    10 = @{ D = "0.07422"; E = "  -4.62%  " }
    11 = @{ D = "1.091"; E = "  -4.38%  " }
    12 = @{ D = "1.009"; E = "  +0.77%  " }
    13 = @{ D = "20.75"; E = "  -7.45%  " }
    14 = @{ D = "6.084"; E = "  -4.58%  " }
    15 = @{ D = "7.323"; E = "  -2.87%  " }
    16 = @{ D = "1.786.23"; E = "  -2.92%  " }
    17 = @{ D = "91.46"; E = "  -2.13%  " }
    18 = @{ D = "0.00001057"; E = "  -2.84%  " }
    19 = @{ D = "0.06397"; E = "  -0.16%  " }
    20 = @{ D = "1.004"; E = "  +0.28%  " }
    21 = @{ D = "17.05"; E = "  -2.99%  " }
    22 = @{ D = "5.963"; E = "  -6.40%  " }
    23 = @{ D = "27.854.98"; E = "  -2.01%  " }
    24 = @{ D = "11.30"; E = "  -3.83%  " }
    25 = @{ D = "2.088"; E = "  -8.12%  " }
    26 = @{ D = "157.73"; E = "  +2.10%  " }
    27 = @{ D = "20.20"; E = "  -3.18%  " }
    28 = @{ D = "1.986.62"; E = "  -2.83%  " }
    29 = @{ D = "2.157"; E = "  -9.37%  " }
    30 = @{ D = "125.29"; E = "  -2.89%  " }
    31 = @{ D = "1.132"; E = "  -6.89%  " }
    32 = @{ D = "3.660"; E = "  -0.31%  " }
    33 = @{ D = "5.585"; E = "  -5.88%  " }
    34 = @{ D = "0.08880"; E = "  -4.56%  " }
    35 = @{ D = "12.37"; E = "  -5.67%  " }
    36 = @{ D = "0.02291"; E = "  -3.16%  " }
    37 = @{ D = "0.2105"; E = "  -4.51%  " }
    38 = @{ D = "4.986"; E = "  -4.58%  " }
    39 = @{ D = "0.06001"; E = "  -4.88%  " }
    40 = @{ D = "0.6314"; E = "  -5.49%  " }
    41 = @{ D = "1.177"; E = "  -1.85%  " }
    42 = @{ D = "1.005"; E = "  +0.38%  " }
    43 = @{ D = "1.401"; E = "  -0.11%  " }
    44 = @{ D = "7.816"; E = "  -4.66%  " }
    45 = @{ D = "13.50"; E = "  -4.00%  " }
    46 = @{ D = "0.5901"; E = "  -4.40%  " }
    47 = @{ D = "3.693"; E = "  -2.41%  " }
    48 = @{ D = "1.991"; E = "  -3.58%  " }
    49 = @{ D = "122.20"; E = "  -4.24%  " }
    50 = @{ D = "1.179"; E = "  +2.31%  " }
    51 = @{ D = "0.06846"; E = "  -2.35%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
